# Generate Report for Handback
#
# This script brings the "localization-status.xlsx" workbook up to date
# after a handback: the Overview/zh-cn/de-de sheets move from
# "Ready for handoff" to "Handed back: in sync with en-US", the per-locale
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated for the handback that just completed, and the
# column widths for the now-longer text are widened to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$sourceMdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b4e3827364ce9c981f11b8806f249d56124f70f/e2e/7cc8ff44-6ac0-4958-9a22-582003046462.md"
$sourceMdDisplay = "7cc8ff44-6ac0-4958-9a22-582003046462.md"

$zhXlf = "7cc8ff44-6ac0-4958-9a22-582003046462.5a34a9afc934c77e49db6a6e5385674c5835d76f.zh-cn.xlf"
$deXlf = "7cc8ff44-6ac0-4958-9a22-582003046462.5a34a9afc934c77e49db6a6e5385674c5835d76f.de-de.xlf"

$zhHandbackTime = "2016-08-26 22:57:58"
$deHandbackTime = "2016-08-26 22:58:09"

$hyperlinkColor = 15570276  # RGB(100,149,237) == #6495ED, the workbook's HyperLink style colour
$xlUnderlineStyleSingle = 2

# ---------------------------------------------------------------------
# 1. Update the "Status" text everywhere it appears (Overview E2/F2,
#    zh-cn C2, de-de C2) now that the handback is complete.
# ---------------------------------------------------------------------
if ($overview.Range("E2").Value2 -eq $oldStatus) { $overview.Range("E2").Value2 = $newStatus }
if ($overview.Range("F2").Value2 -eq $oldStatus) { $overview.Range("F2").Value2 = $newStatus }
if ($zhcn.Range("C2").Value2 -eq $oldStatus)      { $zhcn.Range("C2").Value2 = $newStatus }
if ($dede.Range("C2").Value2 -eq $oldStatus)      { $dede.Range("C2").Value2 = $newStatus }

# ---------------------------------------------------------------------
# 2. zh-cn row: Latest Target File (I2), Latest Handback File (J2),
#    Latest Handback DateTime (K2).
# ---------------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceMdUrl, "", "", $sourceMdDisplay) | Out-Null
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("I2").Font.Underline = $xlUnderlineStyleSingle

$zhcn.Range("J2").Value2 = $zhXlf
$zhcn.Range("K2").Value2 = $zhHandbackTime

# ---------------------------------------------------------------------
# 3. de-de row: Latest Target File (I2), Latest Handback File (J2),
#    Latest Handback DateTime (K2).
# ---------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), $sourceMdUrl, "", "", $sourceMdDisplay) | Out-Null
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("I2").Font.Underline = $xlUnderlineStyleSingle

$dede.Range("J2").Value2 = $deXlf
$dede.Range("K2").Value2 = $deHandbackTime

# ---------------------------------------------------------------------
# 4. Widen the columns that now hold the longer "Handed back: in sync
#    with en-US" status text and the longer target/handback file names,
#    so the new values are fully visible.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.14   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 29.14   # F: de-de status

$zhcn.Columns.Item(3).ColumnWidth = 29.14   # C: Status
$zhcn.Columns.Item(9).ColumnWidth = 39.17   # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.17  # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth = 29.14   # C: Status
$dede.Columns.Item(9).ColumnWidth = 39.17   # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.17  # J: Latest Handback File

Write-Output "Handback report generated."
